# "Add front page & Add "E"s to titles"
# The front page (first slide) already exists in this deck; the
# remaining visible change is that every slide title of the form
# "Exercise <number>" becomes "Exercise E<number>".

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    if ($s.Shapes.HasTitle) {
        $title = $s.Shapes.Title
        if ($title.HasTextFrame -and $title.TextFrame.HasText) {
            $t = $title.TextFrame.TextRange.Text
            if ($t -match "^Exercise\s+(?!E)(\d+)$") {
                $title.TextFrame.TextRange.Text = "Exercise E" + $matches[1]
            }
        }
    }
}
